# Rename task sheets: replace the first space with an underscore and
# remove any remaining spaces (e.g. "Task1 Relaxing Music1" -> "Task1_RelaxingMusic1")

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($oldName -match '^(Task\d+) (.*)$') {
        $prefix = $Matches[1]
        $rest = $Matches[2] -replace ' ', ''
        $newName = "$prefix`_$rest"
        if ($newName -ne $oldName) {
            $ws.Name = $newName
        }
    }
}
